$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 3: "zrušiť administration_view - DONE"
#   -> remove the _GoBack bookmarkStart/bookmarkEnd pair that wraps "DONE".
#      (text/formatting stays identical, only the bookmark markers go away)
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$xml3 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>zrušiť administration_view</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t>DONE</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 4: "Pridať novú rolu – potencionálny člen, vyhodiť ostatné role, ..."
#   -> highlight the whole paragraph yellow (incl. paragraph mark) and append
#      two new runs: " -" and bold " DONE".
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Pridať novú rolu – potencionálny člen</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>, vyhodiť ostatné role, ostáva iba člen a pot.člen</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> DONE</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml4) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 5: "dorobiť funkciu get_user_id, ktorá vráti id-čko prihláseného ..."
#   -> split the single run in two ("...id-č" / "ko...") and drop an empty
#      _GoBack bookmarkStart/bookmarkEnd pair in between (Word's "last edit"
#      marker landing where the cursor was left).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$xml5 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>dorobiť funkciu get_user_id, ktorá vráti id-č</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ko prihláseného usera alebo null ak nie je nikto prihlásený</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($xml5) | Out-Null

Write-Output "done"
